$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos price/volume snapshot (Price column D, Volume(1h) column E).
# Some Price values are plain decimals (e.g. "214.79") that Excel would otherwise
# auto-convert to a number; briefly forcing a text NumberFormat before the write
# keeps them as text, then resetting the style back to Normal avoids leaving any
# stray cell formatting behind.
$ws.Range("D2").Value = "27.109.58"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "1.677.54"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.517"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.255"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.42"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.15%  "
$ws.Range("E10").Value = "  +0.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0887"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").Value = "1.915.39"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "1.685.36"
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("E14").Value = "  +1.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.535"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.77%  "
$ws.Range("D17").Value = "27.117.49"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "238.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").Value = "0.0₃0740"
$ws.Range("E20").Value = "  +1.09%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.86%  "
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0498"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("D32").Value = "1.566.03"
$ws.Range("E32").Value = "  +5.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.46%  "
$ws.Range("E35").Value = "  +0.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.599"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.28%  "
$ws.Range("E37").Value = "  -1.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.929"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0173"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.06%  "
$ws.Range("E40").Value = "  +2.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "68.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.88%  "
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.45%  "
$ws.Range("E44").Value = "  -2.38%  "
$ws.Range("D45").Value = "1.822.74"
$ws.Range("E45").Value = "  +0.36%  "
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("E48").Value = "  +2.92%  "
$ws.Range("E49").Value = "  +1.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.57%  "
$ws.Range("E51").Value = "  +1.80%  "

Write-Host "Updated cryptos list"
